$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.617.28"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "2.527.75"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.566"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.54%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "2.915.16"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.11%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.534.48"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("E17").Value = "  -4.38%  "
$ws.Range("D18").Value = "42.610.10"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").Value = "0.0₃0947"
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.28%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -6.39%  "
$ws.Range("E28").Value = "  -2.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.96%  "
$ws.Range("E33").Value = "  +10.10%  "
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0781"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.98%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.88%  "
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.110"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.117"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.33%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0299"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "2.002.00"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").Value = "2.770.68"
$ws.Range("E49").Value = "  -4.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "79.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.16%  "
